$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns are treated as text so that
# numeric-looking strings (e.g. "1.000", "0.9990") keep their exact formatting
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '25.543.56'
$ws.Range("E2").Value = '  +2.78%  '
$ws.Range("D3").Value = '1.669.29'
$ws.Range("E3").Value = '  +2.24%  '
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '236.86'
$ws.Range("E5").Value = '  +1.11%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.4731'
$ws.Range("E7").Value = '  +0.46%  '
$ws.Range("D8").Value = '0.2594'
$ws.Range("E8").Value = '  +1.67%  '
$ws.Range("D9").Value = '0.06160'
$ws.Range("E9").Value = '  +1.55%  '
$ws.Range("D10").Value = '1.667.78'
$ws.Range("E10").Value = '  +2.12%  '
$ws.Range("D11").Value = '0.07003'
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("D12").Value = '14.77'
$ws.Range("E12").Value = '  +1.52%  '
$ws.Range("D13").Value = '0.5854'
$ws.Range("E13").Value = '  -2.70%  '
$ws.Range("D14").Value = '4.361'
$ws.Range("E14").Value = '  +1.19%  '
$ws.Range("E15").Value = '  +3.60%  '
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").Value = '25.545.60'
$ws.Range("E18").Value = '  +2.78%  '
$ws.Range("D19").Value = '0.000006724'
$ws.Range("E19").Value = '  +2.74%  '
$ws.Range("D20").Value = '11.40'
$ws.Range("E20").Value = '  +2.83%  '
$ws.Range("D21").Value = '1.882.41'
$ws.Range("E21").Value = '  +1.94%  '
$ws.Range("E22").Value = '  +2.38%  '
$ws.Range("D23").Value = '8.762'
$ws.Range("E23").Value = '  +2.69%  '
$ws.Range("D24").Value = '5.225'
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("D25").Value = '137.04'
$ws.Range("E25").Value = '  +3.07%  '
$ws.Range("D26").Value = '14.97'
$ws.Range("E26").Value = '  +1.70%  '
$ws.Range("E27").Value = '  +1.15%  '
$ws.Range("D28").Value = '1.716'
$ws.Range("E28").Value = '  +5.83%  '
$ws.Range("D29").Value = '104.41'
$ws.Range("E29").Value = '  +1.30%  '
$ws.Range("D30").Value = '4.005'
$ws.Range("E30").Value = '  +6.33%  '
$ws.Range("D31").Value = '0.07826'
$ws.Range("E31").Value = '  +1.57%  '
$ws.Range("D32").Value = '3.625'
$ws.Range("E32").Value = '  +3.13%  '
$ws.Range("D33").Value = '0.04296'
$ws.Range("E33").Value = '  +0.68%  '
$ws.Range("E34").Value = '  +1.52%  '
$ws.Range("D35").Value = '0.9534'
$ws.Range("E35").Value = '  +3.96%  '
$ws.Range("D36").Value = '0.6050'
$ws.Range("E36").Value = '  +5.02%  '
$ws.Range("D37").Value = '0.9318'
$ws.Range("E37").Value = '  +15.47%  '
$ws.Range("D38").Value = '2.516'
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.01479'
$ws.Range("E40").Value = '  -3.56%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '1.848'
$ws.Range("E41").Value = '  +5.09%  '
$ws.Range("D42").Value = '99.61'
$ws.Range("E42").Value = '  +2.62%  '
$ws.Range("D43").Value = '0.3744'
$ws.Range("E43").Value = '  +1.97%  '
$ws.Range("D44").Value = '4.891'
$ws.Range("E44").Value = '  +4.21%  '
$ws.Range("D45").Value = '0.1114'
$ws.Range("E45").Value = '  +2.64%  '
$ws.Range("D46").Value = '6.197'
$ws.Range("E46").Value = '  +3.52%  '
$ws.Range("D47").Value = '0.05263'
$ws.Range("E47").Value = '  +1.08%  '
$ws.Range("D48").Value = '29.85'
$ws.Range("E48").Value = '  +1.78%  '
$ws.Range("D49").Value = '7.461'
$ws.Range("E49").Value = '  +4.81%  '
$ws.Range("E50").Value = '  +0.18%  '
$ws.Range("D51").Value = '1.201'
$ws.Range("E51").Value = '  +2.50%  '
